$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.463.35"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.688.50"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.105"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.369"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.86%  "
$ws.Range("D13").Value = "3.168.73"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "63.311.82"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("D17").Value = "2.695.23"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.48%  "
$ws.Range("D28").Value = "0.0₃0858"
$ws.Range("E28").Value = "  -5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "341.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.952"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.623"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0564"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0973"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0242"
$ws.Range("D51").Style = "Normal"
